# R1 release of Module 3 lessons
# Update the cached "last modified" date field text from 6/17/16 to 6/28/16
# across the slide master and all slide layouts.

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "6/17/16") {
                $tr.Text = "6/28/16"
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $lay = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholders $lay.Shapes
}
